$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a scratch cell as a text formula, then copy/paste-special
# as values-only into the target cell. This guarantees the target keeps a plain text
# (inline/shared string) cell -- matching the source data -- instead of Excel auto-
# coercing number-looking strings (e.g. "0.100") into numeric cells, which would both
# change the cell type and drop meaningful trailing zeros.
function Set-TextValue([string]$cellAddr, [string]$val) {
    $scratch = $ws.Range("Z1")
    $scratch.Formula = '="' + $val + '"'
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
    $scratch.Clear()
}

Set-TextValue "D2" '58.582.60'
$ws.Range("E2").Value = '  -2.04%  '

Set-TextValue "D3" '2.626.05'
$ws.Range("E3").Value = '  +0.45%  '

$ws.Range("E4").Value = '  +0.07%  '

Set-TextValue "D5" '534.74'
$ws.Range("E5").Value = '  -0.96%  '

Set-TextValue "D6" '142.41'
$ws.Range("E6").Value = '  -0.13%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("E8").Value = '  -0.03%  '

Set-TextValue "D9" '2.634.22'
$ws.Range("E9").Value = '  +0.39%  '

Set-TextValue "D10" '7.03'
$ws.Range("E10").Value = '  +8.55%  '

Set-TextValue "D11" '0.100'
$ws.Range("E11").Value = '  -2.35%  '

Set-TextValue "D12" '0.334'
$ws.Range("E12").Value = '  -0.50%  '

$ws.Range("E13").Value = '  +0.83%  '

Set-TextValue "D14" '3.090.23'
$ws.Range("E14").Value = '  +0.82%  '

Set-TextValue "D15" '58.521.86'
$ws.Range("E15").Value = '  -1.97%  '

Set-TextValue "D16" '20.86'
$ws.Range("E16").Value = '  +0.51%  '

Set-TextValue "D17" '2.629.58'
$ws.Range("E17").Value = '  +0.94%  '

Set-TextValue "D18" '0.0000132'
$ws.Range("E18").Value = '  -1.26%  '

Set-TextValue "D19" '4.38'
$ws.Range("E19").Value = '  +0.16%  '

Set-TextValue "D20" '335.03'
$ws.Range("E20").Value = '  -2.34%  '

Set-TextValue "D21" '10.16'
$ws.Range("E21").Value = '  -0.09%  '

Set-TextValue "D22" '6.22'
$ws.Range("E22").Value = '  -3.00%  '

$ws.Range("E23").Value = '  -0.18%  '

Set-TextValue "D24" '66.48'
$ws.Range("E24").Value = '  -1.92%  '

Set-TextValue "D25" '0.414'
$ws.Range("E25").Value = '  +0.74%  '

$ws.Range("E26").Value = '  -0.86%  '

Set-TextValue "D27" '0.999'
$ws.Range("E27").Value = '  +0.13%  '

Set-TextValue "D28" '7.12'
$ws.Range("E28").Value = '  -1.94%  '

Set-TextValue "D29" '0.0₃0738'
$ws.Range("E29").Value = '  -1.89%  '

$ws.Range("E30").Value = '  -0.02%  '

$ws.Range("E31").Value = '  -2.23%  '

$ws.Range("E32").Value = '  +0.04%  '

Set-TextValue "D33" '18.76'
$ws.Range("E33").Value = '  -0.94%  '

Set-TextValue "D34" '150.56'
$ws.Range("E34").Value = '  +0.49%  '

Set-TextValue "D35" '3.89'
$ws.Range("E35").Value = '  -2.63%  '

Set-TextValue "D36" '37.11'
$ws.Range("E36").Value = '  -0.31%  '

$ws.Range("E37").Value = '  -1.47%  '

Set-TextValue "D38" '0.825'
$ws.Range("E38").Value = '  -3.38%  '

Set-TextValue "D39" '1.42'
$ws.Range("E39").Value = '  -3.62%  '

Set-TextValue "D40" '0.812'
$ws.Range("E40").Value = '  -2.57%  '

Set-TextValue "D41" '3.58'
$ws.Range("E41").Value = '  +0.66%  '

Set-TextValue "D42" '281.05'
$ws.Range("E42").Value = '  +2.24%  '

$ws.Range("E43").Value = '  +0.03%  '

$ws.Range("E44").Value = '  +0.34%  '

Set-TextValue "D45" '10.68'
$ws.Range("E45").Value = '  -0.54%  '

Set-TextValue "D46" '19.06'
$ws.Range("E46").Value = '  +2.26%  '

Set-TextValue "D47" '0.0531'
$ws.Range("E47").Value = '  +1.03%  '

$ws.Range("E48").Value = '  -2.34%  '

Set-TextValue "D49" '0.0224'
$ws.Range("E49").Value = '  +0.21%  '

Set-TextValue "D50" '1.938.45'
$ws.Range("E50").Value = '  -1.22%  '

Set-TextValue "D51" '4.45'
$ws.Range("E51").Value = '  -2.08%  '
